# Commit: "Switching to backend supplying mp3 files"
# - Anime sheet: G column location formula switches from forward-slash web
#   paths to backslash local paths with a trailing ".mp3" extension.
# - Anime sheet: I column INSERT statement fixes the "arist" typo to
#   "artist" and naturally reflects the new G-column location text.
# - Anime sheet becomes the active / selected tab & cell, replacing
#   "Video Games" as the previously-selected tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Anime")

# --- G2 / G3:G41 : location formula -----------------------------------
$ws.Range("G2").Formula = '="music\"&LOWER(F2)&"\"&A2&".mp3"'
$ws.Range("G3:G41").Formula = '="music\"&LOWER(F3)&"\"&A3&".mp3"'

# --- I2 / I3:I41 : INSERT statement formula ("arist" -> "artist") ------
$ws.Range("I2").Formula = '="INSERT INTO songs (name, property, difficulty, song_name, artist, category, location, video_link) VALUES (''"&A2&"''"&", "&"''"&B2&"''"&", "&"''"&C2&"''"&", "&"''"&D2&"''"&", "&"''"&E2&"''"&", "&"''"&F2&"''"&", "&"''"&G2&"'', "&"''"&H2&"'');"'
$ws.Range("I3:I41").Formula = '="INSERT INTO songs (name, property, difficulty, song_name, artist, category, location, video_link) VALUES (''"&A3&"''"&", "&"''"&B3&"''"&", "&"''"&C3&"''"&", "&"''"&D3&"''"&", "&"''"&E3&"''"&", "&"''"&F3&"''"&", "&"''"&G3&"'', "&"''"&H3&"'');"'

# --- View state: make Anime the active sheet / tab and move the
#     selection to I37 (as in the saved workbook). ----------------------
$ws.Activate()
$ws.Range("I37").Select()
